$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row of "standard tag names" at row 1 without shifting
# existing data down (the rest of the sheet already starts at row 2).
$ws.Range("C1").Value = "頭條要聞"
$ws.Range("D1").Value = "社會"
$ws.Range("E1").Value = "生活"
$ws.Range("F1").Value = "財經"
$ws.Range("G1").Value = "國際"
$ws.Range("H1").Value = "兩岸"
$ws.Range("I1").Value = "娛樂名人"
$ws.Range("J1").Value = "體育"
$ws.Range("K1").Value = "地方"
$ws.Range("L1").Value = "科技"
$ws.Range("M1").Value = "軍事"
$ws.Range("N1").Value = "政治"
$ws.Range("Q1").Value = "娛樂名人"

# Update the active selection to reflect the author's cursor position at save time.
$ws.Range("K17").Select()
